# ----------------------------------------------------------------------
# CompStat weekly report refresh: bump the report volume/date header and
# load the newly collected crime-complaint figures for the 33rd Precinct.
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume/issue number and the week-covering date range ---
$ws.Range("A8").Value = "Volume 30   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/16/2023  Through  1/22/2023"

# --- Crime-complaint figures (rows 14-29, cols C:N) ---
# Cells whose content type / number format changes (e.g. a blank "N/A" cell
# now has a real count, or a count reverts to "N/A") need their NumberFormat
# set explicitly before the value is written so Excel stores the correct
# cell style.

$ws.Range("D14").NumberFormat = '#,##0'
$ws.Range("D14").Value = 1
$ws.Range("E14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E14").Value = -100
$ws.Range("G14").NumberFormat = '#,##0'
$ws.Range("G14").Value = 1
$ws.Range("H14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H14").Value = -100
$ws.Range("J14").NumberFormat = '#,##0'
$ws.Range("J14").Value = 1
$ws.Range("K14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K14").Value = -100
$ws.Range("M23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M23").Value = -100
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("C27").Value = 2
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 1
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E27").Value = 100
$ws.Range("I27").NumberFormat = '#,##0'
$ws.Range("I27").Value = 2
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("D28").Value = 1
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E28").Value = -100
$ws.Range("G28").NumberFormat = '#,##0'
$ws.Range("G28").Value = 1
$ws.Range("H28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H28").Value = -100
$ws.Range("J28").NumberFormat = '#,##0'
$ws.Range("J28").Value = 1
$ws.Range("K28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K28").Value = -100
$ws.Range("D29").NumberFormat = '#,##0'
$ws.Range("D29").Value = 1
$ws.Range("E29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E29").Value = -100
$ws.Range("G29").NumberFormat = '#,##0'
$ws.Range("G29").Value = 1
$ws.Range("H29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H29").Value = -100
$ws.Range("J29").NumberFormat = '#,##0'
$ws.Range("J29").Value = 1
$ws.Range("K29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K29").Value = -100

# Plain value updates (style/format unchanged).
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 12.5
$ws.Range("I16").Value = 8
$ws.Range("J16").Value = 5
$ws.Range("K16").Value = 60
$ws.Range("L16").Value = -27.272727272727
$ws.Range("M16").Value = -57.894736842105
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 133.333333333333
$ws.Range("I17").Value = 18
$ws.Range("J17").Value = 14
$ws.Range("K17").Value = 28.571428571428
$ws.Range("L17").Value = 50
$ws.Range("M17").Value = 157.142857142857
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 12
$ws.Range("H18").Value = 9.090909090909
$ws.Range("I18").Value = 10
$ws.Range("J18").Value = 10
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = -44.444444444444
$ws.Range("M18").Value = 25
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -14.285714285714
$ws.Range("F19").Value = 18
$ws.Range("G19").Value = 24
$ws.Range("H19").Value = -25
$ws.Range("I19").Value = 15
$ws.Range("J19").Value = 15
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = -16.666666666666
$ws.Range("M19").Value = 36.363636363636
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = -44.444444444444
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = -23.529411764705
$ws.Range("I20").Value = 11
$ws.Range("J20").Value = 17
$ws.Range("K20").Value = -35.294117647058
$ws.Range("L20").Value = 120
$ws.Range("M20").Value = 175
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 71
$ws.Range("G21").Value = 79
$ws.Range("H21").Value = -10.126582278481
$ws.Range("I21").Value = 62
$ws.Range("J21").Value = 62
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = -3.125
$ws.Range("M21").Value = 24
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -75
$ws.Range("J22").Value = 2
$ws.Range("D23").Value = 1
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = -80
$ws.Range("J23").Value = 4
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = -47.619047619047
$ws.Range("G24").Value = 113
$ws.Range("H24").Value = -37.168141592920
$ws.Range("I24").Value = 41
$ws.Range("J24").Value = 86
$ws.Range("K24").Value = -52.325581395348
$ws.Range("L24").Value = -4.651162790697
$ws.Range("M24").Value = 64
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -22.222222222222
$ws.Range("F25").Value = 31
$ws.Range("G25").Value = 33
$ws.Range("H25").Value = -6.060606060606
$ws.Range("I25").Value = 22
$ws.Range("J25").Value = 22
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 46.666666666666
$ws.Range("M25").Value = -18.518518518518
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 2
$ws.Range("K27").Value = 0
